# #21-02 - done the query for authorizations
#
# On the "authorizations query" slide, the two rectangles that sit over the
# "Provider ID" / blank value box at the top of the query-builder diagram
# were nudged/resized slightly. PowerPoint records that kind of tweak as a
# brand new shape appended to the end of the spTree (with a new shape id /
# creationId) plus removal of the old shape, rather than an in-place resize,
# so we reproduce it the same way: duplicate the existing shape (to keep its
# exact p:style / text formatting), delete the original, then move the
# duplicate to the new position/size and rename it to match.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(11)

# EMU -> point helper for Shape.Left/Top/Width/Height (COM uses points,
# 1 pt = 12700 EMU); nudge by a hair so float truncation in the host can't
# round the converted EMU value down by one.
function ToPt($emu) {
    return ($emu / 12700.0) + 0.00005
}

# --- "Provider ID" rectangle (was Rectangle 10 @ 5631870,1200724) ---
$old1 = $s.Shapes.Item("Rectangle 10")
$new1 = $old1.Duplicate()
$old1.Delete()
$new1.Name = "Rectangle 75"
$new1.Left = ToPt(5643413)
$new1.Top = ToPt(1163943)
$new1.Width = ToPt(1590968)
$new1.Height = ToPt(375213)

# --- blank value rectangle (was Rectangle 17 @ 7659256,1200725) ---
$old2 = $s.Shapes.Item("Rectangle 17")
$new2 = $old2.Duplicate()
$old2.Delete()
$new2.Name = "Rectangle 76"
$new2.Left = ToPt(7670801)
$new2.Top = ToPt(1167527)
$new2.Width = ToPt(1847190)
$new2.Height = ToPt(359210)
